$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Extend the table by one more year column (L): replicate K3 (year 2020)
# and K4 (6.18) into column L, including their styling.
$ws.Range("L3").Value = 2020
$ws.Range("K3").Copy()
$ws.Range("L3").PasteSpecial(-4122)

$ws.Range("L4").Value = 6.18
$ws.Range("K4").Copy()
$ws.Range("L4").PasteSpecial(-4122)

# Leave the cursor on M12, matching the saved selection in the workbook.
$null = $ws.Range("M12").Select()
